# "add web audio (for touch)"
#
# On the "Questions" slide, the content placeholder's first bullet
# ("Web audio api works?") is split so the trailing "?" becomes its own
# run, and a new sub-bullet ("Yes (25int)") is added directly beneath it
# at the next outline/indent level.

$p = $ppt.ActivePresentation

# Locate the slide + content placeholder holding the "Web audio api works?"
# bullet robustly (rather than hard-coding slide/shape indices).
$targetSlide = $null
$targetShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    for ($j = 1; $j -le $sl.Shapes.Count; $j++) {
        $shp = $sl.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -like "*Web audio*") {
                $targetSlide = $sl
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# The first paragraph reads "Web audio api works?" - find it.
$para1 = $tr.Paragraphs(1)

# Split the trailing "?" off into its own run (matches the diff turning
# "<a:t> works?</a:t>" into "<a:t> works</a:t>" + "<a:t>?</a:t>").
# Paragraph .Text includes a trailing paragraph-mark character, so the
# visible text length is one less than .Length / .Text.Length.
$visibleLen = $para1.Text.Length - 1
$qPos = $para1.Start + $visibleLen - 1
$qMark = $tr.Characters($qPos, 1)
$qMark.Text = "?"

# Re-fetch paragraph 1 (content/position may have shifted) and add the new
# "Yes (25int)" bullet right after it, one indent level deeper.
$para1 = $tr.Paragraphs(1)
$para1.InsertAfter("`rYes (25int)") | Out-Null

$newPara = $tr.Paragraphs(2)
$newPara.IndentLevel = 2

Write-Host "Added 'Yes (25int)' sub-bullet under 'Web audio api works?'"
